$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 17:07"

# Reorder country names (swap/rotate rows) to match the refreshed ranking
$ws.Range("A53").Value = "Portugal"
$ws.Range("A54").Value = "Honduras"
$ws.Range("A72").Value = "Kenia"
$ws.Range("A73").Value = "Afganistan"
$ws.Range("A87").Value = "Birmania"
$ws.Range("A88").Value = "Costa de Marfil"
$ws.Range("A89").Value = "Republica de Macedonia"
$ws.Range("A137").Value = "Sri Lanka"
$ws.Range("A138").Value = "Aruba"
$ws.Range("A215").Value = "Montserrat"
$ws.Range("A216").Value = "Islas Malvinas"

# Update numeric statistics per country row
$ws.Range("B4").Value = 7682785
$ws.Range("C4").Value = 3141
$ws.Range("D4").Value = 4895967
$ws.Range("E4").Value = 2571691
$ws.Range("G4").Value = 95
$ws.Range("H4").Value = 215127
$ws.Range("B5").Value = 6704900
$ws.Range("C5").Value = 22827
$ws.Range("D5").Value = 5678160
$ws.Range("E5").Value = 922953
$ws.Range("G5").Value = 187
$ws.Range("H5").Value = 103787
$ws.Range("B6").Value = 4940706
$ws.Range("C6").Value = 207
$ws.Range("E6").Value = 498613
$ws.Range("G6").Value = 18
$ws.Range("H6").Value = 146791
$ws.Range("B17").Value = 473306
$ws.Range("C17").Value = 1560
$ws.Range("D17").Value = 445418
$ws.Range("E17").Value = 14818
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = 13070
$ws.Range("B18").Value = 387121
$ws.Range("C18").Value = 4172
$ws.Range("D18").Value = 316371
$ws.Range("E18").Value = 61219
$ws.Range("G18").Value = 67
$ws.Range("H18").Value = 9531
$ws.Range("B21").Value = 330263
$ws.Range("C21").Value = 2677
$ws.Range("D21").Value = 234099
$ws.Range("E21").Value = 60134
$ws.Range("G21").Value = 28
$ws.Range("H21").Value = 36030
$ws.Range("B26").Value = 305567
$ws.Range("C26").Value = 910
$ws.Range("E26").Value = 32242
$ws.Range("G26").Value = 9
$ws.Range("H26").Value = 9625
$ws.Range("B29").Value = 169508
$ws.Range("C29").Value = 548
$ws.Range("D29").Value = 142880
$ws.Range("E29").Value = 17117
$ws.Range("G29").Value = 7
$ws.Range("H29").Value = 9511
$ws.Range("B38").Value = 115371
$ws.Range("C38").Value = 317
$ws.Range("D38").Value = 91569
$ws.Range("E38").Value = 21653
$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 2149
$ws.Range("B46").Value = 94870
$ws.Range("C46").Value = 688
$ws.Range("D46").Value = 83480
$ws.Range("E46").Value = 8080
$ws.Range("G46").Value = 8
$ws.Range("H46").Value = 3310
$ws.Range("B53").Value = 80312
$ws.Range("C53").Value = 427
$ws.Range("D53").Value = 50712
$ws.Range("E53").Value = 27568
$ws.Range("G53").Value = 14
$ws.Range("H53").Value = 2032
$ws.Range("B54").Value = 80020
$ws.Range("C54").Value = 391
$ws.Range("D54").Value = 29768
$ws.Range("E54").Value = 47819
$ws.Range("G54").Value = 11
$ws.Range("H54").Value = 2433
$ws.Range("B61").Value = 57732
$ws.Range("C61").Value = 831
$ws.Range("D61").Value = 41938
$ws.Range("E61").Value = 14405
$ws.Range("G61").Value = 14
$ws.Range("H61").Value = 1389
$ws.Range("B72").Value = 39586
$ws.Range("C72").Value = 137
$ws.Range("D72").Value = 27331
$ws.Range("E72").Value = 11512
$ws.Range("G72").Value = 8
$ws.Range("H72").Value = 743
$ws.Range("B73").Value = 39486
$ws.Range("C73").Value = 64
$ws.Range("D73").Value = 32977
$ws.Range("E73").Value = 5042
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = 1467
$ws.Range("B87").Value = 20033
$ws.Range("C87").Value = 1252
$ws.Range("D87").Value = 5782
$ws.Range("E87").Value = 13780
$ws.Range("G87").Value = 27
$ws.Range("H87").Value = 471
$ws.Range("B88").Value = 19885
$ws.Range("D88").Value = 19490
$ws.Range("E88").Value = 275
$ws.Range("H88").Value = 120
$ws.Range("B89").Value = 19096
$ws.Range("C89").Value = 223
$ws.Range("D89").Value = 15645
$ws.Range("E89").Value = 2683
$ws.Range("G89").Value = 8
$ws.Range("H89").Value = 768
$ws.Range("B104").Value = 10789
$ws.Range("C104").Value = 11
$ws.Range("E104").Value = 276
$ws.Range("B116").Value = 7109
$ws.Range("C116").Value = 97
$ws.Range("D116").Value = 2674
$ws.Range("E116").Value = 4312
$ws.Range("G116").Value = 3
$ws.Range("H116").Value = 123
$ws.Range("B119").Value = 5883
$ws.Range("C119").Value = 38
$ws.Range("D119").Value = 5278
$ws.Range("E119").Value = 482
$ws.Range("B133").Value = 4818
$ws.Range("C133").Value = 51
$ws.Range("D133").Value = 2951
$ws.Range("E133").Value = 1784
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = 83
$ws.Range("B137").Value = 4118
$ws.Range("C137").Value = 605
$ws.Range("D137").Value = 3266
$ws.Range("E137").Value = 839
$ws.Range("H137").Value = 13
$ws.Range("B138").Value = 4094
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 3612
$ws.Range("E138").Value = 451
$ws.Range("H138").Value = 31
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
